$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A, B, D, E hold text that looks numeric/date-like ("7111",
# "9/1/2025", "12", "809371823"). Force them to Text format first so the
# assigned values stay plain strings instead of being auto-coerced into
# numbers/dates (matching every other row in the sheet, which stores these
# as text).
$ws.Cells.Item(78, 1).NumberFormat = "@"
$ws.Cells.Item(78, 2).NumberFormat = "@"
$ws.Cells.Item(78, 4).NumberFormat = "@"
$ws.Cells.Item(78, 5).NumberFormat = "@"

$ws.Range("A78").Value = "7111"
$ws.Range("B78").Value = "9/1/2025"
$ws.Range("C78").Value = "VILELA 3699"
$ws.Range("D78").Value = "12"
$ws.Range("E78").Value = "809371823"
$ws.Range("F78").Value = "NEW"
$ws.Range("G78").Value = "Pendiente"
$ws.Range("H78").Value = "Cambiar "
$ws.Range("I78").Value = 1
$ws.Range("J78").Value = "Cambio"
$ws.Range("K78").Value = "Sin equipos"
$ws.Range("L78").Value = "Terminal"
$ws.Range("M78").Value = -58.482817
$ws.Range("N78").Value = -34.550845
$ws.Range("O78").Value = "Saavedra"
$ws.Range("P78").Value = "Capital Norte"
